{"js": "// Fix the typo \"alweays\" -> \"always\" in the tips document\n// (paragraph: \"(  .style.width  doesn't alweays work  )\").\nconst body = context.document.body;\nconst results = body.search(\"alweays\", { matchCase: false, matchWholeWord: true });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"always\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Fix the typo \"alweays\" -> \"always\" in the tips document.\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n# Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n#         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n# Wrap:=1 (wdFindContinue), Replace:=2 (wdReplaceAll)\n$find.Execute(\"alweays\", $false, $true, $false, $false, $false, $true, 1, $false, \"always\", 2)\n"}
